# Updated cryptos list on Sat Sep 16 22:00:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must stay TEXT even when it looks numeric
# (coerces to "@" text format while writing, then restores the default
# "General"/Normal style so the cell ends up exactly like a fresh
# inline/shared string cell with no leftover per-cell style index).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.786.21"
$ws.Range("E2").Value = "  +0.88%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.649.35"
$ws.Range("E3").Value = "  +0.89%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.52%  "

# Row 5 - BNB
Set-TextValue "D5" "216.61"
$ws.Range("E5").Value = "  +1.38%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.34%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.45%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.72%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0627"
$ws.Range("E9").Value = "  +0.17%  "

# Row 10 - Solana
Set-TextValue "D10" "19.27"
$ws.Range("E10").Value = "  +2.18%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0845"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.878.34"
$ws.Range("E12").Value = "  +0.89%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.654.37"
$ws.Range("E13").Value = "  +1.34%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.46%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.52%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.45"
$ws.Range("E16").Value = "  +0.23%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "26.799.58"
$ws.Range("E17").Value = "  +0.91%  "

# Row 18 - ShibaInu
Set-TextValue "D18" "0.0₃0746"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "218.60"
$ws.Range("E19").Value = "  +1.16%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  +0.54%  "

# Row 21 - Toncoin
$ws.Range("E21").Value = "  +13.11%  "

# Row 22 - Uniswap
Set-TextValue "D22" "4.39"
$ws.Range("E22").Value = "  +1.46%  "

# Row 23 - Chainlink
Set-TextValue "D23" "6.33"
$ws.Range("E23").Value = "  +1.03%  "

# Row 24 - Avalanche
$ws.Range("E24").Value = "  +1.35%  "

# Row 25 - Monero
Set-TextValue "D25" "146.02"
$ws.Range("E25").Value = "  -0.75%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.31%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.17%  "

# Row 28 - Cosmos
Set-TextValue "D28" "7.20"
$ws.Range("E28").Value = "  +4.42%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  +1.22%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.0521"
$ws.Range("E30").Value = "  +1.54%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.54%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.37"
$ws.Range("E32").Value = "  +0.54%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "D33" "3.03"
$ws.Range("E33").Value = "  +1.75%  "

# Rows 34 & 35 swap places: LidoDAOToken <-> Maker, with updated values
$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D34" "1.284.22"
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D35" "1.55"
$ws.Range("E35").Value = "  +2.76%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +2.49%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +2.37%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.544"
$ws.Range("E38").Value = "  +6.69%  "

# Row 39 - ARBITRUM
Set-TextValue "D39" "0.835"
$ws.Range("E39").Value = "  +4.66%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.46%  "

# Row 41 - TrustWalletToken
Set-TextValue "D41" "0.819"
$ws.Range("E41").Value = "  +2.50%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -0.86%  "

# Row 43 - FraxShare
Set-TextValue "D43" "5.45"
$ws.Range("E43").Value = "  +1.76%  "

# Row 44 - RocketPoolETH (only price changes, volume(1h) stays the same)
Set-TextValue "D44" "1.789.77"

# Row 45 - Quant
Set-TextValue "D45" "92.27"
$ws.Range("E45").Value = "  -1.09%  "

# Row 46 - Aave
Set-TextValue "D46" "60.08"
$ws.Range("E46").Value = "  +9.20%  "

# Row 47 - RenderToken
Set-TextValue "D47" "1.62"
$ws.Range("E47").Value = "  +0.99%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  +0.16%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +0.62%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "7.82"
$ws.Range("E50").Value = "  +3.22%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0984"
$ws.Range("E51").Value = "  +2.67%  "
